## adjust vom_cost code in data prep
## change the data prep for vom cost of unit relations
##
## Connections table (Table13 on sheet "Connections"):
##   - rename existing "vom_cost" column to "vom_cost_input1"
##   - add a new trailing column "vom_cost_output1"
##
## Units table (Table1 on sheet "Units"):
##   - rename existing "vom_cost" column to "vom_cost_Input1"
##   - add three new trailing columns:
##       "vom_cost_Input2", "vom_cost_Output1", "vom_cost_Output2"

$wb = $excel.ActiveWorkbook

# --- Connections sheet / Table13 -------------------------------------------------
$wsConnections = $wb.Worksheets.Item("Connections")
$tblConnections = $wsConnections.ListObjects.Item(1)

$vomCostConnections = $tblConnections.ListColumns.Item("vom_cost")
$vomCostConnections.Range.Item(1).Value = "vom_cost_input1"

$vomCostOutput1Connections = $tblConnections.ListColumns.Add()
$vomCostOutput1Connections.Range.Item(1).Value = "vom_cost_output1"

# --- Units sheet / Table1 --------------------------------------------------------
$wsUnits = $wb.Worksheets.Item("Units")
$tblUnits = $wsUnits.ListObjects.Item(1)

$vomCostUnits = $tblUnits.ListColumns.Item("vom_cost")
$vomCostUnits.Range.Item(1).Value = "vom_cost_Input1"

$vomCostInput2Units = $tblUnits.ListColumns.Add()
$vomCostInput2Units.Range.Item(1).Value = "vom_cost_Input2"

$vomCostOutput1Units = $tblUnits.ListColumns.Add()
$vomCostOutput1Units.Range.Item(1).Value = "vom_cost_Output1"

$vomCostOutput2Units = $tblUnits.ListColumns.Add()
$vomCostOutput2Units.Range.Item(1).Value = "vom_cost_Output2"
